$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 71: 23-03-2017, Implementer ---
$ws.Range("A62").Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("A71").Value = 42817

$ws.Range("E71").Value = "Implementer"

# Register the new shared strings in the same order the original author
# typed them (Tværkraft before Sikkerhedsfaktor) so the sharedStrings.xml
# index assignment matches.
$ws.Range("F72").Value = "Har lavet implementation af j-unit for Tværkraft klassen"
$ws.Range("F71").Value = "Har lavet implementation af j-unit for Sikkerhedsfaktor klassen"

$ws.Range("G62").Copy()
$ws.Range("G71:H73").PasteSpecial(-4122)
$ws.Range("G75:H75").PasteSpecial(-4122)

$ws.Range("G71").Value = 0.34375
$ws.Range("H71").Value = 0.39583333333333331

# --- Row 72: continuation ---
$ws.Range("G72").Value = 0.39583333333333331
$ws.Range("H72").Value = 0.4236111111111111

# --- Row 73 ---
$ws.Range("E73").Value = "Reviewer"
$ws.Range("F73").Value = "Har prøvet at fikse problem omkring inertimoment som ikke virker"
$ws.Range("G73").Value = 0.51388888888888895
$ws.Range("H73").Value = 0.625

# --- Row 74: daily total ---
$ws.Range("I74").Value = 4.3499999999999996

# --- Row 75: 24-03-2017, Any Role ---
$ws.Range("A62").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("A75").Value = 42818

$ws.Range("E75").Value = "Any Role"
$ws.Range("F75").Value = "Har lavet bugfixing angående mellemregning af Sikkerhedsfaktor"
$ws.Range("G75").Value = 0.34375
$ws.Range("H75").Value = 0.52083333333333337

# --- Row 76: daily total ---
$ws.Range("I76").Value = 4.1500000000000004

# Widen column F to fit the longer activity descriptions.
$ws.Columns.Item(6).ColumnWidth = 58.5

# Scroll / select like the author left the sheet.
$ws.Range("G76").Select()

Write-Host "done"
